$p = $ppt.ActivePresentation

# The deck's single accessible DrawingML theme (ppt/theme/theme1.xml, bound to
# the slide master) is being swapped from the "Integral" / "Red Violet" color
# scheme to the stock "Office Theme" / "Office" color scheme. Font scheme and
# format scheme are identical between the two themes, so only the 12 theme
# colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) need to change.

$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$master = $p.Slides.Item(1).Master
$themeColors = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [math]::Floor($hex / 0x10000) -band 0xFF
    $g = [math]::Floor($hex / 0x100) -band 0xFF
    $b = $hex -band 0xFF
    $themeColors.Item($i).RGB = $r + ($g * 0x100) + ($b * 0x10000)
}
